$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '74.713.41'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '2.834.93'
$ws.Range("E3").Value = '  +9.38%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '188.54'
$ws.Range("E5").Value = '  +1.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '600.10'
$ws.Range("E6").Value = '  +3.37%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.556'
$ws.Range("E8").Value = '  +3.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.193'
$ws.Range("E9").Value = '  -7.46%  '
$ws.Range("D10").Value = '2.832.33'
$ws.Range("E10").Value = '  +9.22%  '
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.371'
$ws.Range("E12").Value = '  +3.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.90'
$ws.Range("E13").Value = '  +1.85%  '
$ws.Range("D14").Value = '3.363.91'
$ws.Range("D15").Value = '74.866.06'
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.29'
$ws.Range("E16").Value = '  +3.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000188'
$ws.Range("E17").Value = '  -2.67%  '
$ws.Range("D18").Value = '2.856.53'
$ws.Range("E18").Value = '  +9.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.14'
$ws.Range("E19").Value = '  +7.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.42'
$ws.Range("E20").Value = '  +6.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.17'
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.27'
$ws.Range("E22").Value = '  -1.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.12'
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.81'
$ws.Range("E26").Value = '  +1.08%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.990.78'
$ws.Range("E27").Value = '  +9.39%  '
$ws.Range("B28").Value = 'NEARProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.22'
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.60'
$ws.Range("E29").Value = '  +4.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000103'
$ws.Range("E30").Value = '  +9.65%  '
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '526.20'
$ws.Range("E32").Value = '  +4.94%  '
$ws.Range("E33").Value = '  +4.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.89'
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.81'
$ws.Range("E35").Value = '  +5.93%  '
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.120'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '20.06'
$ws.Range("E38").Value = '  +4.32%  '
$ws.Range("E39").Value = '  +1.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.28'
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '181.63'
$ws.Range("E42").Value = '  +22.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.08'
$ws.Range("E43").Value = '  +1.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.340'
$ws.Range("E44").Value = '  +6.04%  '
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("E46").Value = '  +7.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.65'
$ws.Range("E47").Value = '  +1.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.36'
$ws.Range("E48").Value = '  -2.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0849'
$ws.Range("E49").Value = '  +4.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.571'
$ws.Range("E50").Value = '  +9.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.75'
$ws.Range("E51").Value = '  +3.30%  '
